$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column widths for Z (26) and AA (27)
$ws.Columns.Item(26).ColumnWidth = 47.77734375
$ws.Columns.Item(27).ColumnWidth = 38.109375

# Z2: "Allocation Method for number of n months for Unit Linked Products"
$ws.Range("Z2").Value = "Allocation Method for number of n months for Unit Linked Products"
$ws.Range("Z2").Font.Size = 7.5
$ws.Range("Z2").Borders.LineStyle = 1
$ws.Range("Z2").Borders.Weight = 2
$ws.Range("Z2").VerticalAlignment = -4108
$ws.Range("Z2").WrapText = $true

# AA2: "Unit Linked Allocation Method for Male and Female Lives"
$ws.Range("AA2").Value = "Unit Linked Allocation Method for Male and Female Lives"
$ws.Range("AA2").Font.Size = 7.5
$ws.Range("AA2").Borders.LineStyle = 1
$ws.Range("AA2").Borders.Weight = 2
$ws.Range("AA2").VerticalAlignment = -4108
$ws.Range("AA2").WrapText = $true

# Z3: "Percentage of Allocation for n months"
$ws.Range("Z3").Value = "Percentage of Allocation for n months"
$ws.Range("Z3").Font.Size = 7.5
$ws.Range("Z3").Borders.LineStyle = 1
$ws.Range("Z3").Borders.Weight = 2
$ws.Range("Z3").VerticalAlignment = -4108
$ws.Range("Z3").WrapText = $true

# AA3: empty cell, same font/alignment but no border
$ws.Range("AA3").Font.Size = 7.5
$ws.Range("AA3").VerticalAlignment = -4108
$ws.Range("AA3").WrapText = $true

# Update view: scroll to show column Q at top-left, select Z2:AA3
$ws.Application.ActiveWindow.ScrollColumn = 17
$ws.Range("Z2:AA3").Select()
